$wb = $excel.ActiveWorkbook

# Work on the "AddValidEmployeeTest" sheet (adding a new TC there)
$ws = $wb.Worksheets.Item("AddValidEmployeeTest")

# Rename header "middlename" -> "Middle Name"
$ws.Range("D1").Value = "Middle Name"

# Change the second data row's first-name value from "jack2" to "john"
$ws.Range("C3").Value = "john"

# Select this sheet and make E3 the active cell/selection, matching the
# author's last edit position
$ws.Activate()
$ws.Range("E3").Select()
